# Task Management System update:
# - Add a new "Auto Reply Sent" column (X)
# - Mark both tasks' status as COMPLETED (column E, the lowercase 'status' column)
# - Record that task T01 got an auto-reply sent ("Yes")
# - Refresh the "Last Updated" timestamp (column W) for both rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header: X1 "Auto Reply Sent" ---
# Copy the formatting of the neighbouring header cell (bold/border/centered)
# so the new header matches the existing header row styling, then set text.
$ws.Range("W1").Copy()
$ws.Range("X1").PasteSpecial(-4122)
$ws.Range("X1").Value = "Auto Reply Sent"

# --- Row 2 (task MOM-20260107-001-T01) ---
$ws.Range("E2").Value = "COMPLETED"
$ws.Range("W2").Value = 46030.80701517361
$ws.Range("X2").Value = "Yes"

# --- Row 3 (task MOM-20260107-001-T02) ---
$ws.Range("E3").Value = "COMPLETED"
$ws.Range("W3").Value = 46030.80679405093
